# DataSource - Emision AP No Enlatada: add CodigoAgente / NUM_GRUPO columns
# and refresh a handful of data values (nueva fila de flota + numero de
# documento / suma asegurada ajustados).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (F:G) before the old "TIPOPOLIZA" column ---
$ws.Columns("F:G").Insert()

# Insert copies formatting from the neighbouring column; the new columns
# should start out blank/unstyled, so clear them before filling in the
# header + the one data point that actually has values.
$ws.Range("F1:G15").Clear()

# Headers
$ws.Range("F1").Value = "CodigoAgente"
$ws.Range("G1").Value = "NUM_GRUPO"

# Only the first data row carries values for the new columns
$ws.Range("F2").Value = 2302
$ws.Range("G2").Value = "Mattioli"

# --- Updated data values ---
$ws.Range("E2").Value = 1785991583
$ws.Range("R2").Value = 550000
$ws.Range("T2").Value = 21840808
$ws.Range("T3").Value = 21840807

# --- Column widths for the two new columns ---
$ws.Columns("F").ColumnWidth = 12.833333333333334
$ws.Columns("G").ColumnWidth = 11.833333333333334

# --- Selection as left by the author ---
$ws.Range("I13").Select()
